# Add a "2021" data column (column M) to the ICT-students table on Sheet1,
# mirroring the existing "2020" column (column L) formatting, then move the
# active selection to P8 (matches the target workbook's sheetView state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone the formatting of column L (rows 2-10, the header/data block)
#        into column M so the new column visually matches the rest of the
#        table (borders, number formats, fonts) before any values go in.
$ws.Range("L2:L10").Copy()
$ws.Range("M2:M10").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Fill in the new "2021" column's values.
$ws.Range("M3").Value  = 2021
$ws.Range("M4").Value  = 952
$ws.Range("M5").Value  = 10437
$ws.Range("M6").Value  = 2253
$ws.Range("M7").Value  = 8184
$ws.Range("M8").Value  = 14020
$ws.Range("M9").Value  = 5139
$ws.Range("M10").Value = 8881

# --- 3. Update the sheet's active cell/selection to P8.
$ws.Range("P8").Select()
